# 4.c.1 worksheet update: extend the year series from 2010-2021 to 2010-2022,
# inserting 2011 and 2012 data points (columns E,F) and appending 2022 (column P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend formatting rightwards by copying the existing "last column" (M)
#        formatting across the newly used columns (through P). This mirrors how
#        the source workbook grew the table and keeps every cellXf shared/deduped.
$ws.Range("M3").Copy()
$ws.Range("N3:P3").PasteSpecial(-4122)

$ws.Range("M4").Copy()
$ws.Range("D4:P4").PasteSpecial(-4122)

$ws.Range("M5").Copy()
$ws.Range("D5:P5").PasteSpecial(-4122)

$ws.Range("M6").Copy()
$ws.Range("D6:P6").PasteSpecial(-4122)

$ws.Range("M7").Copy()
$ws.Range("D7:P7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Row 4 - year headers 2010...2022
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 2011
$ws.Range("F4").Value = 2012
$ws.Range("G4").Value = 2013
$ws.Range("H4").Value = 2014
$ws.Range("I4").Value = 2015
$ws.Range("J4").Value = 2016
$ws.Range("K4").Value = 2017
$ws.Range("L4").Value = 2018
$ws.Range("M4").Value = 2019
$ws.Range("N4").Value = 2020
$ws.Range("O4").Value = 2021
$ws.Range("P4").Value = 2022

# --- 3. Row 5 data
$ws.Range("D5").Value = 87.9
$ws.Range("E5").Value = 89.6
$ws.Range("F5").Value = 87.5
$ws.Range("G5").Value = 88.8
$ws.Range("H5").Value = 89.8
$ws.Range("I5").Value = 94.7
$ws.Range("J5").Value = 91.6
$ws.Range("K5").Value = 93.4
$ws.Range("L5").Value = 93.5
$ws.Range("M5").Value = 93.6
$ws.Range("N5").Value = 94.5
$ws.Range("O5").Value = 93.5
$ws.Range("P5").Value = 94.2

# --- 4. Row 6 data
$ws.Range("D6").Value = 93.6
$ws.Range("E6").Value = 93.3
$ws.Range("F6").Value = 93.9
$ws.Range("G6").Value = 94.3
$ws.Range("H6").Value = 94.4
$ws.Range("I6").Value = 95
$ws.Range("J6").Value = 95.4
$ws.Range("K6").Value = 96
$ws.Range("L6").Value = 96.4
$ws.Range("M6").Value = 96.3
$ws.Range("N6").Value = 96.7
$ws.Range("O6").Value = 96.6
$ws.Range("P6").Value = 96
$ws.Range("P6").NumberFormat = "0.0"

# --- 5. Row 7 data
$ws.Range("D7").Value = 92.9
$ws.Range("E7").Value = 92.8
$ws.Range("F7").Value = 94.1
$ws.Range("G7").Value = 94.8
$ws.Range("H7").Value = 95.3
$ws.Range("I7").Value = 95.9
$ws.Range("J7").Value = 96.9
$ws.Range("K7").Value = 98
$ws.Range("L7").Value = 98
$ws.Range("M7").Value = 98
$ws.Range("N7").Value = 98.2
$ws.Range("O7").Value = 98.1
$ws.Range("P7").Value = 97.5

# --- 6. Row 8 stray cell keeps its value/position (style renumbers automatically
#        because the now-unused old cellXf is dropped on save).
$ws.Range("L8").Value = $ws.Range("L8").Value

# --- 7. Selection / active cell, matching the author's final cursor position.
$ws.Range("Q4").Select()
